$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "73.978.19"
$ws.Range("E2").Value = "  +7.45%  "
$ws.Range("D3").Value = "2.620.75"
$ws.Range("E3").Value = "  +7.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "187.24"
$ws.Range("E5").Value = "  +14.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "582.06"
$ws.Range("E6").Value = "  +3.75%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +4.73%  "
$ws.Range("E9").Value = "  +16.23%  "
$ws.Range("D10").Value = "2.619.91"
$ws.Range("E10").Value = "  +7.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.164"
$ws.Range("E11").Value = "  +1.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.357"
$ws.Range("E12").Value = "  +7.54%  "
$ws.Range("E13").Value = "  +1.87%  "
$ws.Range("E14").Value = "  +5.21%  "
$ws.Range("D15").Value = "73.988.82"
$ws.Range("E15").Value = "  +7.66%  "
$ws.Range("D16").Value = "3.103.07"
$ws.Range("E16").Value = "  +7.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.43"
$ws.Range("E17").Value = "  +12.82%  "
$ws.Range("D18").Value = "2.614.87"
$ws.Range("E18").Value = "  +7.00%  "
$ws.Range("E19").Value = "  +29.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.77"
$ws.Range("E20").Value = "  +10.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "366.90"
$ws.Range("E21").Value = "  +8.25%  "
$ws.Range("E22").Value = "  +18.03%  "
$ws.Range("E23").Value = "  +5.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.98"
$ws.Range("E25").Value = "  +6.90%  "
$ws.Range("E26").Value = "  +9.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.32"
$ws.Range("E27").Value = "  +10.67%  "
$ws.Range("D28").Value = "2.755.88"
$ws.Range("E28").Value = "  +7.28%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").Value = "0.0₃0940"
$ws.Range("E30").Value = "  +13.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "524.71"
$ws.Range("E31").Value = "  +21.18%  "
$ws.Range("E32").Value = "  +15.28%  "
$ws.Range("E33").Value = "  +6.60%  "
$ws.Range("E34").Value = "  +8.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.34"
$ws.Range("E36").Value = "  +1.89%  "
$ws.Range("E37").Value = "  +10.06%  "
$ws.Range("E39").Value = "  +1.41%  "
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.91"
$ws.Range("E41").Value = "  +12.16%  "
$ws.Range("E42").Value = "  +9.78%  "
$ws.Range("E43").Value = "  +7.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "161.82"
$ws.Range("E44").Value = "  +24.58%  "
$ws.Range("E45").Value = "  +13.83%  "
$ws.Range("E46").Value = "  +9.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.90"
$ws.Range("E47").Value = "  +3.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0856"
$ws.Range("E48").Value = "  +18.91%  "
$ws.Range("E49").Value = "  +7.97%  "
$ws.Range("E50").Value = "  +7.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.88"
$ws.Range("E51").Value = "  +22.76%  "
